$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Add New Bank Accounts", "PASSED", "chrome"),
    @("Edit The Bank Accounts", "PASSED", "chrome"),
    @("Delete The Bank Accounts", "PASSED", "chrome"),
    @("", "PASSED", "chrome"),
    @("Add New Bank Accounts", "PASSED", "chrome"),
    @("Edit The Bank Accounts", "PASSED", "chrome"),
    @("Delete The Bank Accounts", "PASSED", "chrome")
)

$startRow = 7
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i

    $cellA = $ws.Cells.Item($row, 1)
    if ($data[$i][0] -eq "") {
        # A plain Value = "" is indistinguishable from clearing the cell, so
        # the empty shared string is produced via the classic "force text"
        # leading-apostrophe entry. That also stamps the cell with a
        # quote-prefix number format, so the style is reset back to Normal
        # afterwards, leaving a plain, unstyled, empty shared-string cell.
        $cellA.Value = "'"
        $cellA.Style = "Normal"
    } else {
        $cellA.Value = $data[$i][0]
    }

    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
